# Add a new "2021年" data row (row 12) to Sheet1, appending after the
# existing "2020年" row (row 11), with the same 50 columns (A:AX).
#
# Columns B, H, I, Q, Z, AG have no data for 2021 (they are blank in the
# source row just like they are blank in earlier rows), but the sheet's
# convention is that every column in the row still has a (blank) cell
# entry, not simply a missing one. Touching a formatting property (here,
# explicitly clearing the (already absent) cell border) is enough to make
# Excel materialize the blank cell without touching the shared style
# table, which keeps the output's styles.xml identical to the original.

$xlNone = -4142        # xlLineStyleNone
$xlPasteFormats = -4122 # xlPasteFormats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A12: row label -------------------------------------------------------
$ws.Range("A12").Value = "2021年"
# Reuse the same look as the other year cells in column A (bold, centered,
# thin box border) by copying the formatting from the previous year's cell.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial($xlPasteFormats)

# --- B12..AX12: data columns -----------------------------------------------
$ws.Range("B12").Borders.LineStyle = $xlNone
$ws.Range("C12").Value = 218
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 357
$ws.Range("F12").Value = 853
$ws.Range("G12").Value = 2062
$ws.Range("H12").Borders.LineStyle = $xlNone
$ws.Range("I12").Borders.LineStyle = $xlNone
$ws.Range("J12").Value = 75950
$ws.Range("K12").Value = 2825
$ws.Range("L12").Value = 4711
$ws.Range("M12").Value = 358
$ws.Range("N12").Value = 184
$ws.Range("O12").Value = 13892
$ws.Range("P12").Value = 1199
$ws.Range("Q12").Borders.LineStyle = $xlNone
$ws.Range("R12").Value = 178
$ws.Range("S12").Value = 2103
$ws.Range("T12").Value = 64
$ws.Range("U12").Value = 1677
$ws.Range("V12").Value = 1022
$ws.Range("W12").Value = 28
$ws.Range("X12").Value = 5676
$ws.Range("Y12").Value = 649
$ws.Range("Z12").Borders.LineStyle = $xlNone
$ws.Range("AA12").Value = 11102
$ws.Range("AB12").Value = 264
$ws.Range("AC12").Value = 59
$ws.Range("AD12").Value = 1567
$ws.Range("AE12").Value = 1545
$ws.Range("AF12").Value = 3
$ws.Range("AG12").Borders.LineStyle = $xlNone
$ws.Range("AH12").Value = 2285
$ws.Range("AI12").Value = 4935
$ws.Range("AJ12").Value = 462
$ws.Range("AK12").Value = 945
$ws.Range("AL12").Value = 39
$ws.Range("AM12").Value = 641
$ws.Range("AN12").Value = 2762
$ws.Range("AO12").Value = 1202
$ws.Range("AP12").Value = 1144
$ws.Range("AQ12").Value = 1924
$ws.Range("AR12").Value = 64
$ws.Range("AS12").Value = 725
$ws.Range("AT12").Value = 3937
$ws.Range("AU12").Value = 147
$ws.Range("AV12").Value = 277
$ws.Range("AW12").Value = 1828
$ws.Range("AX12").Value = 36
